$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 2500
$ws.Range("I74").Value = 2500
$ws.Range("K74").Value = 2500
$ws.Range("M74").Value = -1564

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 2500
$ws.Range("I77").Value = 2500
$ws.Range("K77").Value = 12500
$ws.Range("M77").Value = -7820

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 2598.9
$ws.Range("I107").Value = 2598.9
$ws.Range("K107").Value = 2598.9
$ws.Range("M107").Value = -678.9000000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3578.2
$ws.Range("I2").Value = 858.125
$ws.Range("J2").Value = 6686.857
$ws.Range("K2").Value = 858.125
$ws.Range("L2").Value = 6686.857
$ws.Range("M2").Value = -745.125
$ws.Range("N2").Value = -6912.857

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2357.75
$ws.Range("I45").Value = 1096.9166
$ws.Range("K45").Value = 1096.9166
$ws.Range("M45").Value = -719.9166

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 965.8333
$ws.Range("I61").Value = 965.8333
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 965.8333
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -753.8333
$ws.Range("N61").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1002.75
$ws.Range("I110").Value = 1002.75
$ws.Range("K110").Value = 1002.75
$ws.Range("M110").Value = 1042.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 3578.2
$ws.Range("I116").Value = 858.125
$ws.Range("J116").Value = 6686.857
$ws.Range("K116").Value = 858.125
$ws.Range("L116").Value = 6686.857
$ws.Range("M116").Value = 1435.875
$ws.Range("N116").Value = -11274.857

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 965.8333
$ws.Range("I136").Value = 965.8333
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 2897.4999
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -347.4998999999998
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3578.2
$ws.Range("I3").Value = 858.125
$ws.Range("J3").Value = 6686.857
$ws.Range("K3").Value = 858.125
$ws.Range("L3").Value = 6686.857
$ws.Range("M3").Value = -744.125
$ws.Range("N3").Value = -6914.857

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5890.3335
$ws.Range("I86").Value = 2879
$ws.Range("K86").Value = 2879
$ws.Range("M86").Value = -1756

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 5890.3335
$ws.Range("I89").Value = 2879
$ws.Range("K89").Value = 14395
$ws.Range("M89").Value = -8779

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3529.1765
$ws.Range("I105").Value = 3178.3572
$ws.Range("K105").Value = 3178.3572
$ws.Range("M105").Value = -1431.3572

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1000
$ws.Range("I107").Value = 1000
$ws.Range("K107").Value = 1000
$ws.Range("M107").Value = 920

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H123").Value = 40000
$ws.Range("J123").Value = 40000
$ws.Range("L123").Value = 40000
$ws.Range("N123").Value = -49800

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 929
$ws.Range("I16").Value = 696.6667
$ws.Range("K16").Value = 696.6667
$ws.Range("M16").Value = -409.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 20225.889
$ws.Range("I41").Value = 5059
$ws.Range("K41").Value = 5059
$ws.Range("M41").Value = -4631

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1447.2307
$ws.Range("I105").Value = 1714.8334
$ws.Range("J105").Value = 1217.8572
$ws.Range("K105").Value = 1714.8334
$ws.Range("L105").Value = 1217.8572
$ws.Range("M105").Value = 32.16660000000002
$ws.Range("N105").Value = -4711.8572

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1000
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 1000
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -4840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 929
$ws.Range("I113").Value = 696.6667
$ws.Range("K113").Value = 696.6667
$ws.Range("M113").Value = 1473.3333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 166667230
$ws.Range("I4").Value = 250000460
$ws.Range("J4").Value = 779.75
$ws.Range("K4").Value = 750001380
$ws.Range("L4").Value = 2339.25
$ws.Range("M4").Value = -750001268
$ws.Range("N4").Value = -2563.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H112").Value = 49999
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 49999
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 149997
$ws.Range("M112").ClearContents()
$ws.Range("N112").Value = -152213

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 1033
$ws.Range("J138").Value = 1033
$ws.Range("L138").Value = 3099
$ws.Range("N138").Value = -13379

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 13770.5
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 13770.5
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 13770.5
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -17610.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1656.091
$ws.Range("I113").Value = 1603.5
$ws.Range("K113").Value = 1603.5
$ws.Range("M113").Value = 566.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2619.8
$ws.Range("I132").Value = 2723.5
$ws.Range("J132").Value = 2205
$ws.Range("K132").Value = 8170.5
$ws.Range("L132").Value = 6615
$ws.Range("M132").Value = -5640.5
$ws.Range("N132").Value = -11675

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4200.8
$ws.Range("I136").Value = 4200.8
$ws.Range("K136").Value = 12602.4
$ws.Range("M136").Value = -10052.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3218.8
$ws.Range("I126").Value = 3232
$ws.Range("K126").Value = 9696
$ws.Range("M126").Value = -7226

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H133").Value = 70000
$ws.Range("J133").Value = 70000
$ws.Range("L133").Value = 70000
$ws.Range("N133").Value = -80120

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1304.862
$ws.Range("I136").Value = 1043.15
$ws.Range("J136").Value = 1886.4445
$ws.Range("K136").Value = 3129.45
$ws.Range("L136").Value = 5659.333500000001
$ws.Range("M136").Value = -579.4500000000003
$ws.Range("N136").Value = -10759.3335
